$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent results for the 380 kV case (Case_5_59)
$nRows = 24

$colB = New-Object 'object[,]' $nRows,1
$colB[0,0] = 18.89028610029899
$colB[1,0] = 18.79215134195095
$colB[2,0] = 18.73593723290201
$colB[3,0] = 18.71406519095684
$colB[4,0] = 18.71049644733685
$colB[5,0] = 18.73563804088604
$colB[6,0] = 18.85562267614499
$colB[7,0] = 19.12205066372569
$colB[8,0] = 19.33544772260212
$colB[9,0] = 19.43605230865415
$colB[10,0] = 19.47462864632365
$colB[11,0] = 19.46629966647525
$colB[12,0] = 19.43921655890106
$colB[13,0] = 19.42268898381581
$colB[14,0] = 19.32894162471204
$colB[15,0] = 19.2723149722561
$colB[16,0] = 19.24007873427315
$colB[17,0] = 19.22922228505552
$colB[18,0] = 19.27830861675755
$colB[19,0] = 19.44715873553945
$colB[20,0] = 19.5602938093406
$colB[21,0] = 19.49966639350912
$colB[22,0] = 19.27559789372856
$colB[23,0] = 19.04677360148718
$ws.Range("B2:B25").Value = $colB

$colC = New-Object 'object[,]' $nRows,1
$colC[0,0] = 6.345714281256412
$colC[1,0] = 6.219701808648485
$colC[2,0] = 6.142998056084868
$colC[3,0] = 6.111956571874066
$colC[4,0] = 6.10681660216491
$colC[5,0] = 6.142578482079656
$colC[6,0] = 6.302155016003737
$colC[7,0] = 6.618316540046903
$colC[8,0] = 6.849874373355015
$colC[9,0] = 6.954518858086117
$colC[10,0] = 6.994003238533271
$colC[11,0] = 6.985506550202258
$colC[12,0] = 6.957770348093844
$colC[13,0] = 6.94076137644838
$colC[14,0] = 6.843017979628846
$colC[15,0] = 6.782847173293603
$colC[16,0] = 6.748175762728876
$colC[17,0] = 6.736427178039574
$colC[18,0] = 6.789259269502275
$colC[19,0] = 6.965921331083479
$colC[20,0] = 7.080533147663052
$colC[21,0] = 7.019453654845735
$colC[22,0] = 6.786360601774116
$colC[23,0] = 6.532724494468975
$ws.Range("C2:C25").Value = $colC

$colD = New-Object 'object[,]' $nRows,1
$colD[0,0] = 4.378109198514563
$colD[1,0] = 4.378831251613356
$colD[2,0] = 4.379658938110849
$colD[3,0] = 4.380093312407942
$colD[4,0] = 4.380171316279412
$colD[5,0] = 4.379664402555576
$colD[6,0] = 4.3782786579512
$colD[7,0] = 4.378591483282851
$colD[8,0] = 4.38064264863013
$colD[9,0] = 4.381965514924827
$colD[10,0] = 4.382521985479924
$colD[11,0] = 4.382399678309422
$colD[12,0] = 4.382010186216382
$colD[13,0] = 4.381778827027072
$colD[14,0] = 4.380563995418543
$colD[15,0] = 4.379918201729717
$colD[16,0] = 4.379583478270243
$colD[17,0] = 4.379476468981391
$colD[18,0] = 4.379983150928014
$colD[19,0] = 4.38212308663286
$colD[20,0] = 4.383845078157281
$colD[21,0] = 4.38289660107827
$colD[22,0] = 4.379953673521002
$colD[23,0] = 4.378185079582371
$ws.Range("D2:D25").Value = $colD

$colE = New-Object 'object[,]' $nRows,1
$colE[0,0] = 10.42935597800257
$colE[1,0] = 10.44275322461251
$colE[2,0] = 10.45162027778671
$colE[3,0] = 10.45539529493552
$colE[4,0] = 10.45603190733874
$colE[5,0] = 10.45167053403384
$colE[6,0] = 10.43384253889257
$colE[7,0] = 10.40395032879017
$colE[8,0] = 10.38505352178071
$colE[9,0] = 10.37711716149373
$colE[10,0] = 10.37420634434596
$colE[11,0] = 10.37482904358711
$colE[12,0] = 10.37687579486647
$colE[13,0] = 10.37814178586181
$colE[14,0] = 10.38558542765721
$colE[15,0] = 10.39032060297966
$colE[16,0] = 10.39310628544836
$colE[17,0] = 10.39406015303702
$colE[18,0] = 10.38981010753267
$colE[19,0] = 10.37627205245185
$colE[20,0] = 10.36797485325254
$colE[21,0] = 10.37235296083552
$colE[22,0] = 10.39004070512438
$colE[23,0] = 10.41149692908399
$ws.Range("E2:E25").Value = $colE

$colF = New-Object 'object[,]' $nRows,1
$colF[0,0] = 51.61517047106191
$colF[1,0] = 51.64967152508014
$colF[2,0] = 51.6795364086974
$colF[3,0] = 51.69388899746361
$colF[4,0] = 51.69640401798823
$colF[5,0] = 51.67972113749503
$colF[6,0] = 51.62526453038361
$colF[7,0] = 51.58735985544781
$colF[8,0] = 51.60148158303814
$colF[9,0] = 51.6170006435634
$colF[10,0] = 51.62418216077639
$colF[11,0] = 51.62257751862543
$colF[12,0] = 51.61756534168403
$colF[13,0] = 51.61466505642222
$colF[14,0] = 51.6006501458525
$colF[15,0] = 51.59437995950237
$colF[16,0] = 51.59162981060164
$colF[17,0] = 51.59084580845475
$colF[18,0] = 51.59495882784305
$colF[19,0] = 51.61900215618122
$colF[20,0] = 51.64231919971518
$colF[21,0] = 51.62917995548272
$colF[22,0] = 51.59469445890848
$colF[23,0] = 51.59023916507387
$ws.Range("F2:F25").Value = $colF

$colI = New-Object 'object[,]' $nRows,1
$colI[0,0] = 37.34093666912185
$colI[1,0] = 37.39697471182244
$colI[2,0] = 37.43669907214367
$colI[3,0] = 37.45422327083317
$colI[4,0] = 37.45721382894713
$colI[5,0] = 37.43693000055298
$colI[6,0] = 37.35915468090557
$colI[7,0] = 37.24885406759996
$colI[8,0] = 37.19358929175828
$colI[9,0] = 37.17404894063046
$colI[10,0] = 37.16745477362657
$colI[11,0] = 37.16883912852326
$colI[12,0] = 37.17349029455583
$colI[13,0] = 37.17644414619391
$colI[14,0] = 37.1949789952094
$colI[15,0] = 37.20778390987623
$colI[16,0] = 37.2156760518131
$colI[17,0] = 37.21843871854606
$colI[18,0] = 37.2063662513323
$colI[19,0] = 37.17210227796039
$colI[20,0] = 37.15440282571172
$colI[21,0] = 37.1634198756331
$colI[22,0] = 37.20700552243982
$colI[23,0] = 37.27416923844825
$ws.Range("I2:I25").Value = $colI

$colJ = New-Object 'object[,]' $nRows,1
$colJ[0,0] = 9.982248629945566
$colJ[1,0] = 9.994377578097199
$colJ[2,0] = 10.00223420293594
$colJ[3,0] = 10.00553910023722
$colJ[4,0] = 10.00609412166712
$colJ[5,0] = 10.00227835543138
$colJ[6,0] = 9.986345907708863
$colJ[7,0] = 9.95833665697142
$colJ[8,0] = 9.939710114260828
$colJ[9,0] = 9.931656033431727
$colJ[10,0] = 9.928666127576401
$colJ[11,0] = 9.929307393299412
$colJ[12,0] = 9.93140885121468
$colJ[13,0] = 9.932703860200629
$colJ[14,0] = 9.940244878922043
$colJ[15,0] = 9.94497822049469
$colJ[16,0] = 9.947740189851578
$colJ[17,0] = 9.948682134051989
$colJ[18,0] = 9.944470264604417
$colJ[19,0] = 9.93078997596486
$colJ[20,0] = 9.922198710643139
$colJ[21,0] = 9.926752136152684
$colJ[22,0] = 9.944699784636624
$colJ[23,0] = 9.96556971233607
$ws.Range("J2:J25").Value = $colJ

$colK = New-Object 'object[,]' $nRows,1
$colK[0,0] = 16.93002746922933
$colK[1,0] = 16.86433294254153
$colK[2,0] = 16.82754616970398
$colK[3,0] = 16.81345958079615
$colK[4,0] = 16.81117545139081
$colK[5,0] = 16.82735251690924
$colK[6,0] = 16.90664652862556
$colK[7,0] = 17.08974813503568
$colK[8,0] = 17.24028403889915
$colK[9,0] = 17.31204920798094
$colK[10,0] = 17.33967956459631
$colK[11,0] = 17.3337089564437
$colK[12,0] = 17.31431337225074
$colK[13,0] = 17.30249165099562
$colK[14,0] = 17.23565860850132
$colK[15,0] = 17.19548765523191
$colK[16,0] = 17.17269277518716
$colK[17,0] = 17.16502866665992
$colK[18,0] = 17.19973192570809
$colK[19,0] = 17.31999813940717
$colK[20,0] = 17.40123903334476
$colK[21,0] = 17.35764382910235
$colK[22,0] = 17.19781215536948
$colK[23,0] = 17.03733942315086
$ws.Range("K2:K25").Value = $colK

$colL = New-Object 'object[,]' $nRows,1
$colL[0,0] = 11.80697136171489
$colL[1,0] = 11.82281063360845
$colL[2,0] = 11.83390621630215
$colL[3,0] = 11.83877274431724
$colL[4,0] = 11.83960167633509
$colL[5,0] = 11.83397045059701
$colL[6,0] = 11.81214862646519
$colL[7,0] = 11.7802084241278
$colL[8,0] = 11.76333084699548
$colL[9,0] = 11.75707723805228
$colL[10,0] = 11.75491335339765
$colL[11,0] = 11.75537031116328
$colL[12,0] = 11.75689512439092
$colL[13,0] = 11.75785569495311
$colL[14,0] = 11.76376814415633
$colL[15,0] = 11.7677595947122
$colL[16,0] = 11.77018944882185
$colL[17,0] = 11.77103520010492
$colL[18,0] = 11.7673208254383
$colL[19,0] = 11.7564417119112
$colL[20,0] = 11.75052165636787
$colL[21,0] = 11.7535726035488
$colL[22,0] = 11.76751877210998
$colL[23,0] = 11.78768994376026
$ws.Range("L2:L25").Value = $colL

$colM = New-Object 'object[,]' $nRows,1
$colM[0,0] = 18.25840781137189
$colM[1,0] = 18.26322245136708
$colM[2,0] = 18.26899284731092
$colM[3,0] = 18.27205337791758
$colM[4,0] = 18.27260443792631
$colM[5,0] = 18.26903125011813
$colM[6,0] = 18.25948471358965
$colM[7,0] = 18.26302487900172
$colM[8,0] = 18.27909436816181
$colM[9,0] = 18.2893028902166
$colM[10,0] = 18.29358269226977
$colM[11,0] = 18.29264259120601
$colM[12,0] = 18.2896467103374
$colM[13,0] = 18.28786548380919
$colM[14,0] = 18.27848531293008
$colM[15,0] = 18.27347137202448
$colM[16,0] = 18.2708604790353
$colM[17,0] = 18.27002344168431
$colM[18,0] = 18.27397687984118
$colM[19,0] = 18.29051545880896
$colM[20,0] = 18.3037363718911
$colM[21,0] = 18.29646039071251
$colM[22,0] = 18.27374749325321
$colM[23,0] = 18.25969360619342
$ws.Range("M2:M25").Value = $colM

$colN = New-Object 'object[,]' $nRows,1
$colN[0,0] = 25.60714598775348
$colN[1,0] = 25.66099979129036
$colN[2,0] = 25.69598810386558
$colN[3,0] = 25.71073014396528
$colN[4,0] = 25.71320730430465
$colN[5,0] = 25.69618495942849
$colN[6,0] = 25.625316349617
$colN[7,0] = 25.50155787286496
$colN[8,0] = 25.41985938350758
$colN[9,0] = 25.38468607534722
$colN[10,0] = 25.37165256438008
$colN[11,0] = 25.37444685918932
$colN[12,0] = 25.38360807465348
$colN[13,0] = 25.38925679461297
$colN[14,0] = 25.42219807323722
$colN[15,0] = 25.44291621532018
$colN[16,0] = 25.45502024692745
$colN[17,0] = 25.45915068526267
$colN[18,0] = 25.44069133176791
$colN[19,0] = 25.38090944982922
$colN[20,0] = 25.34350448604333
$colN[21,0] = 25.36331596268069
$colN[22,0] = 25.44169660066263
$colN[23,0] = 25.53341405326477
$ws.Range("N2:N25").Value = $colN

